$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Insert a new column before AV (column 48). Everything from AV..BA
# shifts one place to the right (AV->AW, AW->AX, ... AZ->BA).
# ---------------------------------------------------------------------
$ws.Range("AV1").EntireColumn.Insert()

# Give the brand-new AV11 cell the same number format / style as its
# row-neighbours (AS11/AT11, style used throughout that row block)
# before we put a value in it, so it matches the surrounding column.
$ws.Range("AT11").Copy()
$ws.Range("AV11").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Populate the new "SF_PLAY_MODE" column with its header, example value
# and actual data value (mirrors the neighbouring AV_FILE column).
# ---------------------------------------------------------------------
$ws.Range("AV1").Value = "SF_PLAY_MODE"
$ws.Range("AV2").Value = "menu"
$ws.Range("AV11").Value = "continuous"

Write-Host ("Dimension: " + $ws.UsedRange.Address())
